$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.89"
$ws.Range("E2").Value = "'1.82%"
$ws.Range("D3").Value = "'35.52"
$ws.Range("E3").Value = "'-1.88%"
$ws.Range("D4").Value = "'5.092"
$ws.Range("E4").Value = "'1.21%"
$ws.Range("D5").Value = "'0.08163"
$ws.Range("E5").Value = "'3.88%"
$ws.Range("D6").Value = "'2.054"
$ws.Range("E6").Value = "'-3.32%"
$ws.Range("D7").Value = "'7.946"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.130"
$ws.Range("E8").Value = "'-0.27%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.944"
$ws.Range("E9").Value = "'10.84%"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9255"
$ws.Range("E10").Value = "'0.50%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1061"
$ws.Range("E11").Value = "'11.36%"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1920"
$ws.Range("E12").Value = "'3.90%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09302"
$ws.Range("E13").Value = "'5.74%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03676"
$ws.Range("E14").Value = "'1.94%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09885"
$ws.Range("E15").Value = "'-0.26%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001431"
$ws.Range("E16").Value = "'0.15%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005688"
$ws.Range("E17").Value = "'-0.14%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.477"
$ws.Range("E18").Value = "'0.26%"
$ws.Range("D19").Value = "'0.3393"
$ws.Range("E19").Value = "'0.60%"
$ws.Range("D20").Value = "'0.1301"
$ws.Range("E20").Value = "'-2.71%"
$ws.Range("D21").Value = "'5.092"
$ws.Range("E21").Value = "'-1.55%"
$ws.Range("E22").Value = "'-1.64%"
$ws.Range("D23").Value = "'0.04554"
$ws.Range("E23").Value = "'-0.27%"
$ws.Range("E24").Value = "'-0.55%"
$ws.Range("D25").Value = "'0.004781"
$ws.Range("E25").Value = "'-0.11%"
$ws.Range("D26").Value = "'0.0001250"
$ws.Range("E26").Value = "'-3.96%"
$ws.Range("D27").Value = "'0.0004449"
$ws.Range("E27").Value = "'-6.39%"
$ws.Range("D39").Value = "'0.01963"
$ws.Range("E39").Value = "'5.91%"
$ws.Range("D40").Value = "'0.04882"
$ws.Range("E40").Value = "'3.52%"
$ws.Range("D41").Value = "'0.007556"
$ws.Range("E41").Value = "'-2.91%"
$ws.Range("B42").Value = "Dexo"
$ws.Range("C42").Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo"
$ws.Range("D42").Value = "'0.009900"
$ws.Range("E42").Value = "'28.06%"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1381"
$ws.Range("E43").Value = "'-0.10%"
$ws.Range("D44").Value = "'0.002220"
$ws.Range("E44").Value = "'0.42%"
$ws.Range("D45").Value = "'0.01162"
$ws.Range("E45").Value = "'3.96%"
$ws.Range("D46").Value = "'0.00006616"
$ws.Range("E46").Value = "'3.91%"
$ws.Range("E47").Value = "'-0.05%"
$ws.Range("D48").Value = "'60.02"
$ws.Range("E48").Value = "'15.92%"
$ws.Range("D49").Value = "'0.001500"
$ws.Range("E49").Value = "'-21.09%"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E51").Value = "'-0.05%"
